$d = $word.ActiveDocument

# --- Build replacement XML fragments (exact OOXML matching the target revision) ---
$fechaXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve">Fecha: </w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve">{{ </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t>fecha_actual</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve"> }}</w:t></w:r></w:p>
'@

$yoXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve">Yo </w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve">{{ </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t>nombre_completo</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve"> }} </w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve">identificado(a), con cedula de ciudadanía </w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve">{{ </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t>cedula_ciudadania</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve"> }} </w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t>firmo el presente documento en constancia de conocerlo, de haberlo comprendido y de aceptar las responsabilidades y competencias relacionadas en él.</w:t></w:r></w:p>
'@

$firmaXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:wps="http://schemas.microsoft.com/office/word/2010/wordprocessingShape" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
      <w:r>
        <w:rPr>
          <w:noProof/>
        </w:rPr>
        <mc:AlternateContent>
          <mc:Choice Requires="wps">
            <w:drawing>
              <wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251659264" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="37B5AC2A" wp14:editId="733E2664">
                <wp:simplePos x="0" y="0"/>
                <wp:positionH relativeFrom="column">
                  <wp:posOffset>53340</wp:posOffset>
                </wp:positionH>
                <wp:positionV relativeFrom="paragraph">
                  <wp:posOffset>61595</wp:posOffset>
                </wp:positionV>
                <wp:extent cx="2562225" cy="1266825"/>
                <wp:effectExtent l="0" t="0" r="9525" b="9525"/>
                <wp:wrapNone/>
                <wp:docPr id="993329692" name="Cuadro de texto 1"/>
                <wp:cNvGraphicFramePr/>
                <a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main">
                  <a:graphicData uri="http://schemas.microsoft.com/office/word/2010/wordprocessingShape">
                    <wps:wsp>
                      <wps:cNvSpPr txBox="1"/>
                      <wps:spPr>
                        <a:xfrm>
                          <a:off x="0" y="0"/>
                          <a:ext cx="2562225" cy="1266825"/>
                        </a:xfrm>
                        <a:prstGeom prst="rect">
                          <a:avLst/>
                        </a:prstGeom>
                        <a:solidFill>
                          <a:schemeClr val="lt1"/>
                        </a:solidFill>
                        <a:ln w="6350">
                          <a:noFill/>
                        </a:ln>
                      </wps:spPr>
                      <wps:txbx>
                        <w:txbxContent>
                          <w:p>
                            <w:pPr>
                              <w:jc w:val="center"/>
                            </w:pPr>
                            <w:r>
                              <w:t>{{ firma }}</w:t>
                            </w:r>
                          </w:p>
                        </w:txbxContent>
                      </wps:txbx>
                      <wps:bodyPr rot="0" spcFirstLastPara="0" vertOverflow="overflow" horzOverflow="overflow" vert="horz" wrap="square" lIns="91440" tIns="45720" rIns="91440" bIns="45720" numCol="1" spcCol="0" rtlCol="0" fromWordArt="0" anchor="t" anchorCtr="0" forceAA="0" compatLnSpc="1">
                        <a:prstTxWarp prst="textNoShape">
                          <a:avLst/>
                        </a:prstTxWarp>
                        <a:noAutofit/>
                      </wps:bodyPr>
                    </wps:wsp>
                  </a:graphicData>
                </a:graphic>
              </wp:anchor>
            </w:drawing>
          </mc:Choice>
          <mc:Fallback>
            <w:pict>
              <v:shapetype w14:anchorId="37B5AC2A" id="_x0000_t202" coordsize="21600,21600" o:spt="202" path="m,l,21600r21600,l21600,xe">
                <v:stroke joinstyle="miter"/>
                <v:path gradientshapeok="t" o:connecttype="rect"/>
              </v:shapetype>
              <v:shape id="Cuadro de texto 1" o:spid="_x0000_s1026" type="#_x0000_t202" style="position:absolute;margin-left:4.2pt;margin-top:4.85pt;width:201.75pt;height:99.75pt;z-index:251659264;visibility:visible;mso-wrap-style:square;mso-wrap-distance-left:9pt;mso-wrap-distance-top:0;mso-wrap-distance-right:9pt;mso-wrap-distance-bottom:0;mso-position-horizontal:absolute;mso-position-horizontal-relative:text;mso-position-vertical:absolute;mso-position-vertical-relative:text;v-text-anchor:top" o:gfxdata="UEsDBBQABgAIAAAAIQC2gziS/gAAAOEBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbJSRQU7DMBBF&#10;90jcwfIWJU67QAgl6YK0S0CoHGBkTxKLZGx5TGhvj5O2G0SRWNoz/78nu9wcxkFMGNg6quQqL6RA&#10;0s5Y6ir5vt9lD1JwBDIwOMJKHpHlpr69KfdHjyxSmriSfYz+USnWPY7AufNIadK6MEJMx9ApD/oD&#10;OlTrorhX2lFEilmcO2RdNtjC5xDF9pCuTyYBB5bi6bQ4syoJ3g9WQ0ymaiLzg5KdCXlKLjvcW893&#10;SUOqXwnz5DrgnHtJTxOsQfEKIT7DmDSUCaxw7Rqn8787ZsmRM9e2VmPeBN4uqYvTtW7jvijg9N/y&#10;JsXecLq0q+WD6m8AAAD//wMAUEsDBBQABgAIAAAAIQA4/SH/1gAAAJQBAAALAAAAX3JlbHMvLnJl&#10;bHOkkMFqwzAMhu+DvYPRfXGawxijTi+j0GvpHsDYimMaW0Yy2fr2M4PBMnrbUb/Q94l/f/hMi1qR&#10;JVI2sOt6UJgd+ZiDgffL8ekFlFSbvV0oo4EbChzGx4f9GRdb25HMsYhqlCwG5lrLq9biZkxWOiqY&#10;22YiTra2kYMu1l1tQD30/bPm3wwYN0x18gb45AdQl1tp5j/sFB2T0FQ7R0nTNEV3j6o9feQzro1i&#10;OWA14Fm+Q8a1a8+Bvu/d/dMb2JY5uiPbhG/ktn4cqGU/er3pcvwCAAD//wMAUEsDBBQABgAIAAAA&#10;IQASxoM/LQIAAFUEAAAOAAAAZHJzL2Uyb0RvYy54bWysVE2P2yAQvVfqf0DcGztukm6tOKs0q1SV&#10;ot2VstWeCYbYEmYokNjpr++AnY9ue6p6wTPM8Jh58/D8vmsUOQrratAFHY9SSoTmUNZ6X9DvL+sP&#10;d5Q4z3TJFGhR0JNw9H7x/t28NbnIoAJVCksQRLu8NQWtvDd5kjheiYa5ERihMSjBNsyja/dJaVmL&#10;6I1KsjSdJS3Y0ljgwjncfeiDdBHxpRTcP0nphCeqoFibj6uN6y6syWLO8r1lpqr5UAb7hyoaVmu8&#10;9AL1wDwjB1v/AdXU3IID6UccmgSkrLmIPWA34/RNN9uKGRF7QXKcudDk/h8sfzxuzbMlvvsCHQ4w&#10;ENIalzvcDP100jbhi5USjCOFpwttovOE42Y2nWVZNqWEY2yczWZ36CBOcj1urPNfBTQkGAW1OJdI&#10;FztunO9TzynhNgeqLte1UtEJWhArZcmR4RSVj0Ui+G9ZSpO2oLOP0zQCawjHe2SlsZZrU8Hy3a4b&#10;Ot1BeUICLPTacIavayxyw5x/ZhbFgD2jwP0TLlIBXgKDRUkF9uff9kM+zgijlLQoroK6HwdmBSXq&#10;m8bpfR5PJkGN0ZlMP2Xo2NvI7jaiD80KsPMxPiXDoxnyvTqb0kLziu9gGW7FENMc7y6oP5sr30se&#10;3xEXy2VMQv0Z5jd6a3iADkyHEbx0r8yaYU4eR/wIZxmy/M24+txwUsPy4EHWcZaB4J7VgXfUblTD&#10;8M7C47j1Y9b1b7D4BQAA//8DAFBLAwQUAAYACAAAACEAMfwNz98AAAAHAQAADwAAAGRycy9kb3du&#10;cmV2LnhtbEyOTU+DQBRF9yb+h8kzcWPsAK22II/GGLWJO4sfcTdlnkBk3hBmCvjvHVe6vLk35558&#10;O5tOjDS41jJCvIhAEFdWt1wjvJQPlxsQzivWqrNMCN/kYFucnuQq03biZxr3vhYBwi5TCI33fSal&#10;qxoyyi1sTxy6TzsY5UMcaqkHNQW46WQSRdfSqJbDQ6N6umuo+tofDcLHRf3+5ObH12l5tezvd2O5&#10;ftMl4vnZfHsDwtPs/8bwqx/UoQhOB3tk7USHsFmFIUK6BhHaVRynIA4ISZQmIItc/vcvfgAAAP//&#10;AwBQSwECLQAUAAYACAAAACEAtoM4kv4AAADhAQAAEwAAAAAAAAAAAAAAAAAAAAAAW0NvbnRlbnRf&#10;VHlwZXNdLnhtbFBLAQItABQABgAIAAAAIQA4/SH/1gAAAJQBAAALAAAAAAAAAAAAAAAAAC8BAABf&#10;cmVscy8ucmVsc1BLAQItABQABgAIAAAAIQASxoM/LQIAAFUEAAAOAAAAAAAAAAAAAAAAAC4CAABk&#10;cnMvZTJvRG9jLnhtbFBLAQItABQABgAIAAAAIQAx/A3P3wAAAAcBAAAPAAAAAAAAAAAAAAAAAIcE&#10;AABkcnMvZG93bnJldi54bWxQSwUGAAAAAAQABADzAAAAkwUAAAAA&#10;" fillcolor="white [3201]" stroked="f" strokeweight=".5pt">
                <v:textbox>
                  <w:txbxContent>
                    <w:p>
                      <w:pPr>
                        <w:jc w:val="center"/>
                      </w:pPr>
                      <w:r>
                        <w:t>{{ firma }}</w:t>
                      </w:r>
                    </w:p>
                  </w:txbxContent>
                </v:textbox>
              </v:shape>
            </w:pict>
          </mc:Fallback>
        </mc:AlternateContent>
      </w:r>
    </w:p>
'@

# --- Locate the 'Fecha: ' paragraph and splice in the {{ fecha_actual }} placeholder ---
$fechaIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.StartsWith("Fecha:")) {
        $fechaIdx = $i
        break
    }
}
if ($fechaIdx -ne -1) {
    $d.Paragraphs.Item($fechaIdx).Range.InsertXML($fechaXml)
}

# --- Locate the 'Yo ___..." paragraph and splice in the template placeholders ---
$yoIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.StartsWith("Yo ___")) {
        $yoIdx = $i
        break
    }
}
if ($yoIdx -ne -1) {
    $d.Paragraphs.Item($yoIdx).Range.InsertXML($yoXml)
}

# --- Replace the trailing empty paragraph with the {{ firma }} signature text box ---
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertXML($firmaXml)
